$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.130.98'
$ws.Range('E2').Value = '  +4.65%  '
$ws.Range('D3').Value = '3.506.87'
$ws.Range('E3').Value = '  +2.37%  '
$ws.Range('E4').Value = '  -0.27%  '
$ws.Range('D5').Value = "'418.45"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.30%  '
$ws.Range('D6').Value = "'132.91"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.98%  '
$ws.Range('D7').Value = "'0.653"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +4.49%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').Value = "'0.781"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +7.38%  '
$ws.Range('D10').Value = "'0.162"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +15.59%  '
$ws.Range('D11').Value = "'43.67"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.26%  '
$ws.Range('D12').Value = "'0.0000266"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +20.50%  '
$ws.Range('D13').Value = "'10.14"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +9.99%  '
$ws.Range('D14').Value = '4.060.59'
$ws.Range('E14').Value = '  +2.20%  '
$ws.Range('E15').Value = '  +0.41%  '
$ws.Range('D16').Value = "'20.61"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.44%  '
$ws.Range('D17').Value = '3.511.83'
$ws.Range('E17').Value = '  +2.74%  '
$ws.Range('D18').Value = "'12.81"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.27%  '
$ws.Range('D19').Value = "'1.12"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.06%  '
$ws.Range('D20').Value = '64.996.15'
$ws.Range('E20').Value = '  +4.43%  '
$ws.Range('D21').Value = "'453.98"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.78%  '
$ws.Range('D22').Value = "'90.39"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.35%  '
$ws.Range('E23').Value = '  -1.30%  '
$ws.Range('D24').Value = "'13.36"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.52%  '
$ws.Range('D25').Value = "'3.44"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +4.75%  '
$ws.Range('D26').Value = "'10.00"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.29%  '
$ws.Range('D27').Value = "'34.20"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.40%  '
$ws.Range('D28').Value = "'12.66"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +6.85%  '
$ws.Range('E29').Value = '  +3.11%  '
$ws.Range('E30').Value = '  -1.57%  '
$ws.Range('E31').Value = '  +5.36%  '
$ws.Range('E32').Value = '  -2.22%  '
$ws.Range('D33').Value = "'39.99"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.77%  '
$ws.Range('D34').Value = "'0.999"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.05%  '
$ws.Range('D35').Value = "'57.15"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.47%  '
$ws.Range('D36').Value = "'0.0511"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +4.88%  '
$ws.Range('D37').Value = '0.0₃0744'
$ws.Range('E37').Value = '  +38.90%  '
$ws.Range('E38').Value = '  +11.52%  '
$ws.Range('E39').Value = '  -0.33%  '
$ws.Range('E40').Value = '  +0.21%  '
$ws.Range('D41').Value = "'4.56"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +6.82%  '
$ws.Range('E42').Value = '  +5.51%  '
$ws.Range('D43').Value = "'146.19"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.87%  '
$ws.Range('D44').Value = "'3.32"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.10%  '
$ws.Range('D45').Value = "'0.314"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.33%  '
$ws.Range('E46').Value = '  -3.09%  '
$ws.Range('E47').Value = '  +1.00%  '
$ws.Range('D48').Value = "'15.93"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.49%  '
$ws.Range('D49').Value = "'0.146"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.59%  '
$ws.Range('D50').Value = "'2.57"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +11.46%  '
$ws.Range('D51').Value = "'21.66"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.77%  '
